$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("JLC")

# Insert a new row at 30 (shifts existing rows 30-40 down to 31-41)
$ws.Rows.Item(30).Insert()

# Populate the new part row (R118 - 1k5 0402 pull resistor for USB D+)
# Edited in the same column order the original author used so that new
# shared-string entries land at the same indices as the authoritative file.
$ws.Cells.Item(30, 1).Value2 = 1                       # A30 Qty
$ws.Cells.Item(30, 2).Value2 = "R118"                  # B30 ID
$ws.Cells.Item(30, 3).Value2 = "1k5"                   # C30 VALUE
$ws.Cells.Item(30, 5).Value2 = "Uniroyal Elec"         # E30 MF
$ws.Cells.Item(30, 7).Value2 = "R0403"                 # G30 PACKAGE
$ws.Cells.Item(30, 9).Value2 = "NO"                    # I30 EXTENDED
$ws.Cells.Item(30, 10).Value2 = 0.0009                 # J30 PRICE_1PLUS
$ws.Cells.Item(30, 17).Value2 = 0.0003                 # Q30 PRICE_500PLUS
$ws.Cells.Item(30, 20).Formula = "=A30*J30+ (A30*H30*0.0015)"  # T30 SUMA
$ws.Cells.Item(30, 21).Formula = "=A30*50"             # U30 Pcs for 50 PCB
$ws.Cells.Item(30, 25).Value2 = "CHIP RESISTOR ±1% 1/16W"      # Y30 DESCRIPTION
$ws.Cells.Item(30, 19).Value2 = 43941                  # S30 PRICE_DATE
$ws.Cells.Item(30, 8).Value2 = 2                       # H30 PADS
$ws.Cells.Item(30, 6).Value2 = "C25867"                # F30 LCSC_ID
$ws.Cells.Item(30, 4).Value2 = "0402WGF1501TCE"        # D30 MF_ID
